$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16) with the new built-in
#    table style GUID.
# ---------------------------------------------------------------------------
$newTableStyle = "{9C741C08-136E-4EA8-92BE-96032AAAA147}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the deck's applied theme from "Integral" (Red Violet) back to the
#    default "Office Theme" colour palette. The font/format schemes are
#    already identical between the two themes, so only the 12 theme colour
#    slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) need updating.
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
